# Relaatio taulukko - rename columns to lowercase / snake_case naming.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PELAAJA table (row 3-4)
$ws.Range("B4").Value = "pelaajan_id"
$ws.Range("C4").Value = "nimi"
$ws.Range("D4").Value = "puhnum"
$ws.Range("E4").Value = "kotipaikka"

# RATA table (row 8-9)
$ws.Range("B9").Value = "radan_id"
$ws.Range("C9").Value = "luokitus"
$ws.Range("D9").Value = "vaylien_lkm"
$ws.Range("E9").Value = "osoite"
$ws.Range("F9").Value = "ratamestari"

# PELI table (row 13-14)
$ws.Range("B14").Value = "pelin_id"
$ws.Range("C14").Value = "radan_id"
$ws.Range("D14").Value = "paivamaara"

# PELAAMASSA table (row 18-19)
$ws.Range("B19").Value = "pelin_id"
$ws.Range("C19").Value = "pelaajan_id"

# SUORITUS table (row 23-24)
$ws.Range("B24").Value = "pelaajan_id"
$ws.Range("C24").Value = "pelin_id"
$ws.Range("D24").Value = "radan_id"
$ws.Range("E24").Value = "vaylannumero"

# VAYLA table (row 28-29)
$ws.Range("B29").Value = "radan_id"
$ws.Range("C29").Value = "par"
$ws.Range("D29").Value = "numero"
$ws.Range("E29").Value = "pituus"

# Update the view's active selection to match the saved workbook state.
$ws.Range("H28").Select()
